$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = "27.698.28"
$c.ClearFormats()
$ws.Range("E2").Value = "  +0.19%  "
$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = "1.584.07"
$c.ClearFormats()
$ws.Range("E3").Value = "  -2.06%  "
$ws.Range("E4").Value = "  +1.35%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "206.76"
$c.ClearFormats()
$ws.Range("E5").Value = "  -1.24%  "
$ws.Range("E6").Value = "  -1.89%  "
$ws.Range("E7").Value = "  +1.38%  "
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = "22.24"
$c.ClearFormats()
$ws.Range("E8").Value = "  -3.73%  "
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "0.252"
$c.ClearFormats()
$ws.Range("E9").Value = "  -0.90%  "
$ws.Range("E10").Value = "  -2.43%  "
$ws.Range("E11").Value = "  -0.83%  "
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "1.808.83"
$c.ClearFormats()
$ws.Range("E12").Value = "  -2.05%  "
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = "1.600.69"
$c.ClearFormats()
$ws.Range("E13").Value = "  -1.02%  "
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "3.88"
$c.ClearFormats()
$ws.Range("E14").Value = "  -2.57%  "
$ws.Range("E15").Value = "  -4.32%  "
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "27.657.22"
$c.ClearFormats()
$ws.Range("E16").Value = "  +0.01%  "
$ws.Range("E17").Value = "  -2.05%  "
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "219.01"
$c.ClearFormats()
$ws.Range("E18").Value = "  -3.69%  "
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "0.0₃0693"
$c.ClearFormats()
$ws.Range("E19").Value = "  -3.04%  "
$ws.Range("E20").Value = "  -4.31%  "
$ws.Range("E21").Value = "  +1.34%  "
$ws.Range("E22").Value = "  -3.71%  "
$ws.Range("E23").Value = "  -5.50%  "
$ws.Range("E24").Value = "  -2.81%  "
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "154.79"
$c.ClearFormats()
$ws.Range("E25").Value = "  +0.36%  "
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "6.83"
$c.ClearFormats()
$ws.Range("E26").Value = "  -0.77%  "
$ws.Range("E27").Value = "  +1.34%  "
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "15.10"
$c.ClearFormats()
$ws.Range("E28").Value = "  -1.93%  "
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "0.107"
$c.ClearFormats()
$ws.Range("E29").Value = "  -3.28%  "
$ws.Range("E30").Value = "  -1.71%  "
$ws.Range("E31").Value = "  -2.49%  "
$ws.Range("E32").Value = "  -3.53%  "
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "1.380.28"
$c.ClearFormats()
$ws.Range("E33").Value = "  -0.78%  "
$ws.Range("E34").Value = "  -4.65%  "
$ws.Range("E35").Value = "  -3.95%  "
$ws.Range("E36").Value = "  -2.69%  "
$ws.Range("E37").Value = "  +0.30%  "
$ws.Range("E38").Value = "  -2.81%  "
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "0.536"
$c.ClearFormats()
$ws.Range("E39").Value = "  -3.41%  "
$ws.Range("E40").Value = "  -2.56%  "
$ws.Range("E41").Value = "  +1.39%  "
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "0.975"
$c.ClearFormats()
$ws.Range("E42").Value = "  -3.38%  "
$ws.Range("B43").Value = "RenderToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "1.77"
$c.ClearFormats()
$ws.Range("E43").Value = "  -2.68%  "
$ws.Range("B44").Value = "MXToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "2.17"
$c.ClearFormats()
$ws.Range("E44").Value = "  +2.82%  "
$ws.Range("B45").Value = "Aave"
$ws.Range("C45").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "63.48"
$c.ClearFormats()
$ws.Range("E45").Value = "  -3.03%  "
$ws.Range("E46").Value = "  -2.69%  "
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "1.719.96"
$c.ClearFormats()
$ws.Range("E47").Value = "  -2.04%  "
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "88.31"
$c.ClearFormats()
$ws.Range("E48").Value = "  +0.89%  "
$ws.Range("E49").Value = "  +10.72%  "
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "0.0972"
$c.ClearFormats()
$ws.Range("E50").Value = "  -3.69%  "
$ws.Range("E51").Value = "  -0.54%  "
